$d = $word.ActiveDocument

$pairs = @(
    @("340×2=680", "149×6=894"),
    @("264×3=792", "400×6=2400"),
    @("702×3=2106", "106×7=742"),
    @("767×9=6903", "507×4=2028"),
    @("559×8=4472", "941×5=4705"),
    @("124×8=992", "206×4=824"),
    @("821×3=2463", "215×7=1505"),
    @("546×5=2730", "741×5=3705"),
    @("688×7=4816", "309×2=618"),
    @("233×9=2097", "479×8=3832"),
    @("884×5=4420", "801×8=6408"),
    @("784×2=1568", "563×3=1689"),
    @("512×3=1536", "415×9=3735"),
    @("702×2=1404", "218×7=1526"),
    @("613×7=4291", "231×4=924"),
    @("463×2=926", "715×7=5005"),
    @("488×5=2440", "327×2=654"),
    @("675×6=4050", "999×6=5994"),
    @("141×6=846", "191×7=1337"),
    @("886×3=2658", "695×3=2085"),
    @("413×3=1239", "340×9=3060"),
    @("544×6=3264", "663×6=3978"),
    @("925×4=3700", "535×7=3745"),
    @("317×3=951", "161×4=644"),
    @("295×9=2655", "989×6=5934")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
